$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that was on A2 (restaurant link to a specific chain
# restaurant) and clear the now-stale restaurant link values in A2/A3 — the
# "skip chain restaurants" fix drops the explicit restaurant links.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()

# Replace the old restaurant-name / restaurant-location filter values with
# the new search filters.
$ws.Range("B2").Value = "Japanese"
$ws.Range("B3").Value = "Korean"
$ws.Range("C2").Value = "Kwun Tong"
$ws.Range("C3").Value = "Causeway Bay"

# Turn off the scrape toggles that were left on, and fix the restaurants
# limit (raised from 5 to 10) now that chain-restaurant results are skipped.
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# Update the remembered selection to match the author's last cursor position.
$ws.Range("F10").Select()
